$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from op2 to wong3
$ws.Name = "wong3"

# Update cell values per the diff (experimental results without normalization, Wong3)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 3

$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 51

$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 5
$ws.Range("H5").Value = 8

$ws.Range("H7").Value = 13

$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 2
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 3
$ws.Range("H9").Value = 3

$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 3
$ws.Range("H16").Value = 3

$ws.Range("H21").Value = 10

$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 3
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = 84

$ws.Range("J24").Value = 4
$ws.Range("K24").Value = 84

$ws.Range("B29").Value = 8
$ws.Range("C29").Value = 8
$ws.Range("E29").Value = 8
$ws.Range("F29").Value = 8

$ws.Range("E31").Value = 4

$ws.Range("E33").Value = 4
$ws.Range("F33").Value = 4
$ws.Range("H33").Value = 5

$ws.Range("H34").Value = 9

$ws.Range("J37").Value = 2
$ws.Range("K37").Value = 51

$ws.Range("E40").Value = 6
$ws.Range("F40").Value = 6
$ws.Range("H40").Value = 7

$ws.Range("E46").Value = 4
$ws.Range("F46").Value = 4
$ws.Range("H46").Value = 5

$ws.Range("E48").Value = 2
$ws.Range("F48").Value = 2
$ws.Range("H48").Value = 4
$ws.Range("J48").Value = 3
$ws.Range("K48").Value = 79

$ws.Range("B50").Value = 5
$ws.Range("C50").Value = 5
$ws.Range("E50").Value = 5
$ws.Range("F50").Value = 5
$ws.Range("H50").Value = 12

$ws.Range("E51").Value = 7
$ws.Range("F51").Value = 7
$ws.Range("H51").Value = 8

$ws.Range("E52").Value = 6
$ws.Range("F52").Value = 6
$ws.Range("J52").Value = 3
$ws.Range("K52").Value = 79
